# Updates the cryptocurrency price/volume table (coinranking.com scrape) to
# the latest snapshot. For the "Price" column (D) we prefix the literal with
# an apostrophe and then reset NumberFormat to "General" - this forces Excel
# to store the value as text (matching the source data's inline-string cell
# type) instead of silently reinterpreting e.g. "7.200" / "1.000" as the
# number 7.2 / 1 and dropping the trailing zero.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.489.86"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +3.44%  "

# Row 3
$ws.Range("D3").Value = "'1.816.89"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +4.54%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.35%  "

# Row 5
$ws.Range("D5").Value = "'343.92"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +3.39%  "

# Row 6
$ws.Range("E6").Value = "  +0.56%  "

# Row 7
$ws.Range("D7").Value = "'0.3835"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +2.95%  "

# Row 8
$ws.Range("D8").Value = "'0.3539"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +4.33%  "

# Row 9
$ws.Range("D9").Value = "'49.05"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.13%  "

# Row 10
$ws.Range("D10").Value = "'1.235"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +2.62%  "

# Row 11
$ws.Range("D11").Value = "'0.07777"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +3.54%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.40%  "

# Row 13
$ws.Range("E13").Value = "  +9.06%  "

# Row 14
$ws.Range("D14").Value = "'6.593"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.30%  "

# Row 15
$ws.Range("D15").Value = "'1.815.17"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +4.48%  "

# Row 16
$ws.Range("D16").Value = "'7.200"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +2.12%  "

# Row 17
$ws.Range("E17").Value = "  +2.43%  "

# Row 18
$ws.Range("D18").Value = "'0.06727"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$ws.Range("D19").Value = "'86.52"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.14%  "

# Row 20
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.47%  "

# Row 21
$ws.Range("D21").Value = "'17.66"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +5.67%  "

# Row 22
$ws.Range("E22").Value = "  +5.86%  "

# Row 23
$ws.Range("D23").Value = "'13.15"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("D24").Value = "'27.493.77"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.46%  "

# Row 25
$ws.Range("D25").Value = "'2.467"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.25%  "

# Row 26
$ws.Range("D26").Value = "'2.688"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +7.07%  "

# Row 27
$ws.Range("D27").Value = "'22.19"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +12.67%  "

# Row 28
$ws.Range("D28").Value = "'1.465"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +3.67%  "

# Row 29
$ws.Range("D29").Value = "'154.10"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.54%  "

# Row 30
$ws.Range("D30").Value = "'2.020.96"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +4.69%  "

# Row 31
$ws.Range("D31").Value = "'135.99"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +3.05%  "

# Row 32
$ws.Range("D32").Value = "'6.368"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +2.61%  "

# Row 33
$ws.Range("D33").Value = "'4.075"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.24%  "

# Row 34
$ws.Range("D34").Value = "'13.90"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +5.71%  "

# Row 35
$ws.Range("D35").Value = "'0.08806"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +3.22%  "

# Row 36
$ws.Range("E36").Value = "  -1.60%  "

# Row 37
$ws.Range("D37").Value = "'5.628"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +2.91%  "

# Row 38
$ws.Range("D38").Value = "'0.7048"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +12.67%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06517"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.80%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2259"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +3.81%  "

# Row 41
$ws.Range("D41").Value = "'0.02402"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +2.17%  "

# Row 42
$ws.Range("D42").Value = "'8.984"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +3.90%  "

# Row 43
$ws.Range("D43").Value = "'1.295"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +4.38%  "

# Row 44
$ws.Range("D44").Value = "'14.89"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +2.15%  "

# Row 45
$ws.Range("D45").Value = "'0.6628"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +8.98%  "

# Row 46
$ws.Range("D46").Value = "'1.000"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.44%  "

# Row 47
$ws.Range("D47").Value = "'3.964"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.77%  "

# Row 48
$ws.Range("D48").Value = "'2.189"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +6.30%  "

# Row 49
$ws.Range("D49").Value = "'132.44"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +2.72%  "

# Row 50
$ws.Range("D50").Value = "'0.07330"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("D51").Value = "'80.79"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +3.70%  "
